$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.007.84"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.44%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.484.57"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.27"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.24"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -4.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -4.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.484.94"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -6.60%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.412"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.080.98"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.95"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -6.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.095.75"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.34%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.477.78"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.93"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.86"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "366.99"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -7.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.80"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.538"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.10"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000122"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.72"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -6.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.178"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.64%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "24.24"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.78"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -5.69%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -7.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.08"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.56"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "29.61"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +12.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "158.93"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.887"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.783.93"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.56"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -10.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.42"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.93%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -6.49%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.04%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.24"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -7.93%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "307.31"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -6.40%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.40%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.83%  "
